# Rename the embedded logo pictures' underlying drawing object names.
#
#   * The two BTec logo pictures (one in the "default" header, one in the
#     "first page" header) are currently exported with the drawing name
#     "image2.jpg" -> they should become "image1.jpg".
#   * The two Pearson logo pictures (one in the "default" footer, one in
#     the "first page" footer) are currently exported with the drawing
#     name "image1.png" -> they should become "image2.png".
#
# InlineShape has no writable Name property in the Word object model, so
# each picture is temporarily converted to a floating Shape (which does
# expose Name), renamed, then converted back to an inline shape so the
# surrounding paragraph/run layout is unchanged.

$d = $word.ActiveDocument

function Rename-InlineLogo($range, $newName) {
    $count = $range.InlineShapes.Count
    for ($i = 1; $i -le $count; $i++) {
        $inlineShape = $range.InlineShapes.Item($i)
        $shape = $inlineShape.ConvertToShape()
        $shape.Name = $newName
        [void]$shape.ConvertToInlineShape()
    }
}

foreach ($sec in $d.Sections) {

    # Headers: BTec_Logo-Orange picture, "image2.jpg" -> "image1.jpg"
    for ($hIdx = 1; $hIdx -le 3; $hIdx++) {
        $hdr = $sec.Headers.Item($hIdx)
        if ($hdr.Exists) {
            Rename-InlineLogo $hdr.Range "image1.jpg"
        }
    }

    # Footers: Pearson logo picture, "image1.png" -> "image2.png"
    for ($fIdx = 1; $fIdx -le 3; $fIdx++) {
        $ftr = $sec.Footers.Item($fIdx)
        if ($ftr.Exists) {
            Rename-InlineLogo $ftr.Range "image2.png"
        }
    }
}

Write-Host "Logo picture names updated."
